$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "test"
$ws.Range("A3").Value = "test1"
